$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching the style of the
# existing header row (bold font, thin box border, centered/top aligned).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# I0 / IF data values for rows 2-52.
$data = @(
    @(9,9),
    @(8,8),
    @(8,8),
    @(8,8),
    @(6,6),
    @(8,8),
    @(6,6),
    @(9,9),
    @(6,6),
    @(6,6),
    @(9,9),
    @(5,5),
    @(8,8),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(9,9),
    @(6,7),
    @(6,6),
    @(12,12),
    @(10,10),
    @(5,5),
    @(8,8),
    @(7,7),
    @(6,6),
    @(12,12),
    @(8,8),
    @(6,6),
    @(5,5),
    @(8,8),
    @(7,7),
    @(7,7),
    @(9,9),
    @(8,8),
    @(8,8),
    @(6,6),
    @(8,8),
    @(8,8),
    @(9,9),
    @(7,7),
    @(10,10),
    @(7,8),
    @(7,8),
    @(6,6),
    @(7,8),
    @(6,7),
    @(6,6),
    @(9,9),
    @(7,7),
    @(5,5)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}
